$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quotation row appended at the bottom of the log (row 72).
# Leading apostrophes force Excel to store the date/time as literal
# text instead of auto-converting them to date/time serial numbers,
# matching the existing rows' text-based layout. The style is reset
# afterwards so the new cells don't pick up a "quoted text" format
# that the rest of the sheet doesn't use.
$ws.Range("A72").Value = "'2025-10-11"
$ws.Range("B72").Value = "'21:19:55"
$ws.Range("C72").Value = "1.00 EUR = 1,756.2048"
$ws.Range("A72:C72").Style = "Normal"
